$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'256.49"
$ws.Range('E2').Formula = "'0.08%"
$ws.Range('D3').Formula = "'26.57"
$ws.Range('E3').Formula = "'-1.98%"
$ws.Range('D4').Formula = "'4.689"
$ws.Range('E4').Formula = "'1.36%"
$ws.Range('D5').Formula = "'0.05924"
$ws.Range('E5').Formula = "'0.37%"
$ws.Range('D6').Formula = "'6.608"
$ws.Range('E6').Formula = "'-0.46%"
$ws.Range('D7').Formula = "'0.8543"
$ws.Range('E7').Formula = "'-1.68%"
$ws.Range('D8').Formula = "'0.9138"
$ws.Range('E8').Formula = "'-3.61%"
$ws.Range('D9').Formula = "'0.1378"
$ws.Range('E9').Formula = "'-1.83%"
$ws.Range('D10').Formula = "'0.04345"
$ws.Range('E10').Formula = "'15.86%"
$ws.Range('D11').Formula = "'0.06997"
$ws.Range('E11').Formula = "'-1.20%"
$ws.Range('D12').Formula = "'0.03036"
$ws.Range('E12').Formula = "'-5.10%"
$ws.Range('D13').Formula = "'0.09111"
$ws.Range('E13').Formula = "'-1.63%"
$ws.Range('D14').Formula = "'0.001532"
$ws.Range('E14').Formula = "'-0.99%"
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Formula = "'0.006015"
$ws.Range('E15').Formula = "'-1.90%"
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Formula = "'3.473"
$ws.Range('E16').Formula = "'-1.09%"
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Formula = "'3.138"
$ws.Range('E17').Formula = "'-1.75%"
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').Formula = "'2.168"
$ws.Range('E18').Formula = "'-2.30%"
$ws.Range('B19').Value = 'One'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D19').Formula = "'0.01031"
$ws.Range('E19').Formula = "'1,622.91%"
$ws.Range('D20').Formula = "'0.3081"
$ws.Range('E20').Formula = "'-0.95%"
$ws.Range('D21').Formula = "'0.1284"
$ws.Range('E21').Formula = "'0.05%"
$ws.Range('D22').Formula = "'3.880"
$ws.Range('E22').Formula = "'0.67%"
$ws.Range('D23').Formula = "'0.04207"
$ws.Range('E23').Formula = "'-1.02%"
$ws.Range('D24').Formula = "'0.001213"
$ws.Range('E24').Formula = "'-0.77%"
$ws.Range('D25').Formula = "'0.004655"
$ws.Range('E25').Formula = "'8.66%"
$ws.Range('E26').Formula = "'-0.16%"
$ws.Range('D27').Formula = "'0.0001522"
$ws.Range('E27').Formula = "'1.31%"
$ws.Range('D40').Formula = "'0.03796"
$ws.Range('E40').Formula = "'-0.43%"
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').Formula = "'0.006242"
$ws.Range('E41').Formula = "'0.49%"
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Formula = "'0.1099"
$ws.Range('E42').Formula = "'-0.09%"
$ws.Range('D43').Formula = "'0.002309"
$ws.Range('E43').Formula = "'15.41%"
$ws.Range('D44').Formula = "'0.01454"
$ws.Range('E44').Formula = "'29.68%"
$ws.Range('D45').Formula = "'0.00005128"
$ws.Range('E45').Formula = "'-6.80%"
$ws.Range('D46').Formula = "'0.00000000749"
$ws.Range('E46').Formula = "'-0.15%"
$ws.Range('D47').Formula = "'0.04994"
$ws.Range('E47').Formula = "'-37.98%"
$ws.Range('E48').Formula = "'10,453.78%"
$ws.Range('D49').Formula = "'0.00002098"
$ws.Range('E49').Formula = "'-0.15%"
$ws.Range('D50').Formula = "'0.0001998"
$ws.Range('E50').Formula = "'-0.15%"
